# GO4036_Gizis.docx edit: "added C4 and C5 abstracts"
#
# 1. Merge the run holding the manual line break right before "John " with
#    the run holding "John " itself (they already share identical run
#    formatting, so Word collapses them into a single <w:r>).
# 2. Paragraph 2 (the abstract body) loses its "justify" alignment.
# 3. The (hidden) _GoBack bookmark shifts a little earlier in paragraph 2's
#    text, from right before "at the stellar/" to right before
#    "e characterize variability ..." (splitting "unique" into "uniqu" + "e").

$d = $word.ActiveDocument

# --- Change 1: merge the <w:br/> run with the following "John " run -------
$brk = $d.Content
$brk.Find.Execute("John ") | Out-Null
$johnStart = $brk.Start

# The break character immediately precedes "John "; remove its run boundary
# by deleting the break char and retyping break+"John " as one contiguous
# range so the engine emits a single run.
$breakRange = $d.Range($johnStart - 1, $johnStart)
$breakRange.Delete()

$johnRange = $d.Range($johnStart - 1, $johnStart - 1 + 5)
$johnRange.Text = [char]11 + "John "

# --- Change 2: drop the "both" (justify) alignment on paragraph 2 ---------
$d.Paragraphs.Item(2).Alignment = 0

# --- Change 3: move the _GoBack bookmark earlier in the text --------------
$findRange = $d.Content
$findRange.Find.Execute("K2 can uniqu") | Out-Null
$newBookmarkPos = $findRange.End
$bookmarkRange = $d.Range($newBookmarkPos, $newBookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
